# edit.ps1
# Applies the CasosColombia.xlsx update: refresh several "NaN"/value cells in
# columns DT/CU/BZ/AI (rows 9-195) and append a new data row 202 (date 44096).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Scalar cell fix-ups scattered through the sheet (CU.. columns mostly).
#    Each entry is "CellRef|NewValue". A value of "NaN" is written as the
#    literal text "NaN" (matches the shared string already used elsewhere in
#    the sheet); everything else is a plain number.
# ---------------------------------------------------------------------------
$cellEdits = @(
    "DT9|1",
    "CU17|1",
    "BZ18|NaN",
    "CU26|2",
    "CU27|3",
    "CU28|7",
    "CU29|7",
    "CU30|7",
    "CU31|7",
    "CU32|7",
    "CU33|8",
    "CU34|8",
    "CU35|NaN",
    "CU36|NaN",
    "CU37|NaN",
    "CU38|12",
    "CU39|12",
    "CU40|12",
    "CU41|12",
    "CU42|NaN",
    "CU43|NaN",
    "CU44|18",
    "CU45|19",
    "CU46|19",
    "CU47|19",
    "CU48|20",
    "CU49|22",
    "CU50|22",
    "CU51|24",
    "CU52|24",
    "CU53|24",
    "CU54|24",
    "CU55|24",
    "CU56|24",
    "CU57|24",
    "CU58|25",
    "CU59|25",
    "CU60|25",
    "CU61|25",
    "CU62|25",
    "CU63|26",
    "CU64|26",
    "CU65|26",
    "CU66|26",
    "CU67|26",
    "CU68|27",
    "CU69|28",
    "CU70|28",
    "CU71|32",
    "CU72|32",
    "CU73|32",
    "CU74|32",
    "CU75|32",
    "CU76|35",
    "CU77|36",
    "CU78|36",
    "CU79|36",
    "CU80|39",
    "CU81|40",
    "CU82|40",
    "CU83|43",
    "CU84|45",
    "CU85|47",
    "CU86|48",
    "CU87|50",
    "CU88|53",
    "CU89|57",
    "CU90|62",
    "CU91|69",
    "CU92|70",
    "CU93|72",
    "CU94|73",
    "CU119|141",
    "CU120|145",
    "CU121|147",
    "CU122|154",
    "CU123|158",
    "CU124|161",
    "CU125|164",
    "CU126|171",
    "CU127|182",
    "CU128|185",
    "CU129|188",
    "CU130|189",
    "CU131|195",
    "CU132|202",
    "CU133|210",
    "CU134|213",
    "CU135|213",
    "CU136|217",
    "CU137|220",
    "CU138|223",
    "CU139|240",
    "CU140|245",
    "CU141|264",
    "CU142|276",
    "CU143|280",
    "CU144|288",
    "CU145|291",
    "CU146|306",
    "CU147|313",
    "CU148|331",
    "CU152|391",
    "CU153|402",
    "CU154|418",
    "CU155|438",
    "CU156|449",
    "CU157|454",
    "CU158|474",
    "CU159|477",
    "CU160|485",
    "CU161|505",
    "CU162|513",
    "CU163|530",
    "CU164|552",
    "CU165|569",
    "CU166|574",
    "CU170|640",
    "CU171|647",
    "CU172|661",
    "CU173|676",
    "CU174|704",
    "AI195|NaN"
)

foreach ($entry in $cellEdits) {
    $parts = $entry.Split("|")
    $ref = $parts[0]
    $val = $parts[1]
    if ($val -eq "NaN") {
        $ws.Range($ref).Value = "NaN"
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# ---------------------------------------------------------------------------
# 2) Append the new row 202 (2020-09-06 / serial 44096) at the bottom of the
#    table, columns A:DX, reusing the same values captured from the source.
# ---------------------------------------------------------------------------
$row202 = @(
    "A202|44096",
    "B202|777537",
    "C202|2738",
    "D202|105690",
    "E202|66783",
    "F202|256114",
    "G202|28270",
    "H202|6079",
    "I202|4894",
    "J202|7793",
    "K202|8549",
    "L202|17552",
    "M202|3965",
    "N202|23251",
    "O202|31091",
    "P202|7529",
    "Q202|9664",
    "R202|14692",
    "S202|13661",
    "T202|17634",
    "U202|14791",
    "V202|3668",
    "W202|2804",
    "X202|9700",
    "Y202|28260",
    "Z202|13767",
    "AA202|11161",
    "AB202|58048",
    "AC202|1950",
    "AD202|1014",
    "AE202|712",
    "AF202|469",
    "AG202|651",
    "AH202|455",
    "AI202|631",
    "AJ202|2033",
    "AK202|5182",
    "AL202|37711",
    "AM202|9226",
    "AN202|2541",
    "AO202|45123",
    "AP202|1094",
    "AQ202|22636",
    "AR202|1524",
    "AS202|10067",
    "AT202|1647",
    "AU202|1602",
    "AV202|7681",
    "AW202|1985",
    "AX202|956",
    "AY202|2497",
    "AZ202|2664",
    "BA202|61380",
    "BB202|13893",
    "BC202|5755",
    "BD202|9526",
    "BE202|6512",
    "BF202|277",
    "BG202|1460",
    "BH202|2717",
    "BI202|743",
    "BJ202|2147",
    "BK202|9621",
    "BL202|9462",
    "BM202|10325",
    "BN202|14247",
    "BO202|1963",
    "BP202|898",
    "BQ202|12896",
    "BR202|10604",
    "BS202|12466",
    "BT202|2631",
    "BU202|2082",
    "BV202|5419",
    "BW202|4617",
    "BX202|2026",
    "BY202|5696",
    "BZ202|3426",
    "CA202|2026",
    "CB202|941",
    "CC202|2886",
    "CD202|2211",
    "CE202|1868",
    "CF202|1582",
    "CG202|6039",
    "CH202|2060",
    "CI202|1427",
    "CJ202|1746",
    "CK202|2062",
    "CL202|2096",
    "CM202|2480",
    "CN202|1655",
    "CO202|1210",
    "CP202|1204",
    "CQ202|939",
    "CR202|3372",
    "CS202|1428",
    "CT202|947",
    "CU202|1026",
    "CV202|1710",
    "CW202|1555",
    "CX202|765",
    "CY202|865",
    "CZ202|1278",
    "DA202|1582",
    "DB202|1487",
    "DC202|1532",
    "DD202|1190",
    "DE202|334",
    "DF202|365",
    "DG202|804",
    "DH202|758",
    "DI202|480",
    "DJ202|543",
    "DK202|381",
    "DL202|666",
    "DM202|750",
    "DN202|527",
    "DO202|491",
    "DP202|374",
    "DQ202|521",
    "DR202|134826",
    "DS202|329491",
    "DT202|17671",
    "DU202|142698",
    "DV202|88196",
    "DW202|43100",
    "DX202|12112"
)

foreach ($entry in $row202) {
    $parts = $entry.Split("|")
    $ref = $parts[0]
    $val = $parts[1]
    $ws.Range($ref).Value = [double]$val
}

# ---------------------------------------------------------------------------
# 3) Keep the frozen header pane, move the visible/active cell down to the
#    newly-added row so the saved view matches the source workbook.
# ---------------------------------------------------------------------------
$ws.Range("B202").Select()
